# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders/rewrites the worker-period detail rows (B16:G48) so the two
# workers' monthly periods are interleaved in ascending period order
# instead of being grouped per-worker in descending period order, and
# updates the mora (F) / salary (G) values to match the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Each row: RowNumber, DocType, DocNumber, WorkerName, Period, ValorMora, SalarioBasico
$rows = @(
    @(16, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1704", 14754,  781242),
    @(17, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1705", 14754,  781242),
    @(18, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1706", 14754,  781242),
    @(19, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1707", 14754,  781242),
    @(20, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1708", 14754,  781242),
    @(21, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1709", 14754,  781242),
    @(22, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1710", 14754,  781242),
    @(23, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1711", 14754,  781242),
    @(24, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1712", 14754,  781242),
    @(25, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1801", 22132,  781242),
    @(26, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1802", 22132,  781242),
    @(27, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1803", 22132,  781242),
    @(28, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1804", 22132,  781242),
    @(29, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1805", 22132,  781242),
    @(30, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1806", 22132,  781242),
    @(31, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1807", 22132,  781242),
    @(32, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1808", 22132,  781242),
    @(33, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1809", 23437,  781242),
    @(34, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1810", 23437,  781242),
    @(35, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1811", 23437,  781242),
    @(36, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1812", 23437,  781242),
    @(37, "CC", "45458022",   "ROSA AMELIA CASTRO NARVAEZ", "1901", 80000, 2000000),
    @(38, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1901", 31249,  781242),
    @(39, "CC", "45458022",   "ROSA AMELIA CASTRO NARVAEZ", "1902", 80000, 2000000),
    @(40, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1902", 31249,  781242),
    @(41, "CC", "45458022",   "ROSA AMELIA CASTRO NARVAEZ", "1903", 80000, 2000000),
    @(42, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1903", 31249,  781242),
    @(43, "CC", "45458022",   "ROSA AMELIA CASTRO NARVAEZ", "1904", 80000, 2000000),
    @(44, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1904", 31249,  781242),
    @(45, "CC", "45458022",   "ROSA AMELIA CASTRO NARVAEZ", "1905", 80000, 2000000),
    @(46, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1905", 31249,  781242),
    @(47, "CC", "45458022",   "ROSA AMELIA CASTRO NARVAEZ", "1906", 45334, 2000000),
    @(48, "CC", "1047454774", "DAVID RICARDO ANAYA PEREZ", "1906", 17708,  781242)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]  # B: Tipo Doc Trabajador
    $ws.Cells.Item($rowNum, 3).Value = $r[2]  # C: N Doc Trabajador
    $ws.Cells.Item($rowNum, 4).Value = $r[3]  # D: Nombre Trabajador
    $ws.Cells.Item($rowNum, 5).Value = $r[4]  # E: Periodo Mora
    $ws.Cells.Item($rowNum, 6).Value = $r[5]  # F: Valor Mora
    $ws.Cells.Item($rowNum, 7).Value = $r[6]  # G: Salario Basico
}
